$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.850.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.829.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.453"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("E12").Value = "  -0.59%  "

$ws.Range("E13").Value = "  -0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.466.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.879.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.77%  "

$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.872.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.54%  "

$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("E23").Value = "  -3.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.978.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("E40").Value = "  +0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.04%  "

$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.51%  "

$ws.Range("E46").Value = "  +10.87%  "

$ws.Range("E47").Value = "  +12.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "148.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "389.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
